# daily auto push: 2026-01-19 09:47 UTC
# Insert a new day-row for 2026/01/19 (time slot 16) just above the
# existing "2026/12/29" block, which pushes rows 670:711 down to
# 671:712 (the sheet's used range grows from D711 to D712).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 670 and everything below it down by one row.
$ws.Rows("670:670").Insert()

# Column A holds plain text dates (e.g. "2026/01/19"), not real Excel
# date serials, so force Text format on the new cell before writing it
# - otherwise the date-look-alike string gets auto-converted to a
# date number by the Value setter.
$ws.Range("A670").NumberFormat = "@"
$ws.Range("A670").Value = "2026/01/19"
$ws.Range("B670").Value = "月"
$ws.Range("C670").Value = 16
$ws.Range("D670").Value = 201
